$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# --- Rows 1-4: simple text replacements ---
$table.Cell(1,1).Range.Text = "0M"
$table.Cell(2,1).Range.Text = "0M"
$table.Cell(3,1).Range.Text = "0M"
$table.Cell(4,1).Range.Text = "182"

# --- Insert 3 new rows after row 4 (before old row 5) ---
$newRow1 = $table.Rows.Add($table.Rows.Item(5))
$newRow1.Cells.Item(1).Range.Text = "0.00002"

$newRow2 = $table.Rows.Add($table.Rows.Item(6))
$newRow2.Cells.Item(1).Range.Text = "0.00067"

$newRow3 = $table.Rows.Add($table.Rows.Item(7))
$newRow3.Cells.Item(1).Range.Text = "0.00016"

# Now: row 8 = old row5 (0.00004, unchanged)
#      row 9 = old row6 (0.00037 -> 0.00027)
#      row10 = old row7 (0.00013 -> 0.00034)
#      row11 = old row8 (0.00004 -> 0.00042)
#      row12 = old row9 (0.00016, to be deleted)
#      row13 = old row10 (0.00018, to be deleted)
#      row14 = old row11 (0.00021, to be deleted)
#      row15 = old row12 (0.01123 -> 0.03429)
#      row16 = old row13 (100.0, unchanged)

$table.Cell(9,1).Range.Text = "0.00027"
$table.Cell(10,1).Range.Text = "0.00034"
$table.Cell(11,1).Range.Text = "0.00042"

# Delete the 3 now-obsolete rows (12,13,14 contain 0.00016 / 0.00018 / 0.00021)
$table.Rows.Item(12).Delete()
$table.Rows.Item(12).Delete()
$table.Rows.Item(12).Delete()

# Remaining row (was row15, now row12 after deletions) -> 0.03429
$table.Cell(12,1).Range.Text = "0.03429"

# --- Collapse the trailing multi-value rows (now rows 44-46) into single values ---
$table.Cell(44,1).Range.Text = "99.94"
$table.Cell(45,1).Range.Text = "0.03"
$table.Cell(46,1).Range.Text = "53"
